# Rename the worksheets:
#   "Sheet1"  -> "Data"
#   "KabKota" -> "Rumah Sakit"
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Sheet1")
$wsRS   = $wb.Worksheets.Item("KabKota")

$wsData.Name = "Data"
$wsRS.Name   = "Rumah Sakit"

# Move the selected cell on the "Rumah Sakit" sheet from E11 to B11.
[void]$wsRS.Range("B11").Select()

# Move the selected cell on the "Data" sheet from A2 to F2, and leave "Data"
# as the active (tab-selected) sheet, matching the original workbook state.
[void]$wsData.Select()
[void]$wsData.Range("F2").Select()
